$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Run/Status flags: D2 no->yes, D3 yes->no, D5 no->yes (D4 unchanged)
$ws.Range("D2").Value = "yes"
$ws.Range("D3").Value = "no"
$ws.Range("D5").Value = "yes"

# Move the active selection from D4 to D5
$ws.Range("D5").Select()
